$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The author's edit removes 3 blank spacer rows from the gap between the
# "Выгрузка" block (ends row 38) and the signature block (used to start at
# row 49, now starts at row 46). Deleting entire rows 39:41 shifts every row
# at/after 42 up by 3 (49->46, 50->47, 51->48, 52->49, ... 60->57) and Excel
# re-keys cell/row references, data validations and the used dimension
# automatically.
# ---------------------------------------------------------------------------
$ws.Rows("39:41").Delete()

# The row that lands on 45 (just above the reinstated signature block) keeps
# a short explicit height in the final file.
$ws.Rows("45").RowHeight = 15.75

# The conditional-format rule that used to cover F29:L49 must keep tracking
# the same (now shifted) block -- F29:L46 -- without losing its dxf/priority.
$cf = $ws.Range("F29:L49").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("F29:L46"))

# ---------------------------------------------------------------------------
# Defined names that pointed below row 38 need to be re-pointed at their new
# (shifted) cells; names above the deleted block are untouched by Excel
# already and don't need any action.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    switch ($n.Name) {
        "Port_Letter!Print_Area" { $n.RefersTo = "=Port_Letter!`$A`$1:`$D`$56" }
        "Seal_seller_end"        { $n.RefersTo = "=Port_Letter!`$C`$48" }
        "Seal_seller_start"      { $n.RefersTo = "=Port_Letter!`$C`$46" }
        "Sign_seller_start"      { $n.RefersTo = "=Port_Letter!`$C`$47" }
    }
}
# The Cyrillic-named ranges below the signature block also shifted by 3 rows;
# addressed positionally since some consoles mangle non-ASCII literals.
$wb.Names.Item(10).RefersTo = "=Port_Letter!`$A`$48"   # Письмо_дата
$wb.Names.Item(14).RefersTo = "=Port_Letter!`$D`$47"   # Подписант
$wb.Names.Item(15).RefersTo = "=Port_Letter!`$A`$47"   # Подписант_комментарий

# ---------------------------------------------------------------------------
# View state: selection moves to C44 (was C38) now that the sheet is shorter.
# ---------------------------------------------------------------------------
$ws.Range("C44").Select()

# ---------------------------------------------------------------------------
# Print scale bumped from 85% to 88% to keep the (now shorter) sheet fitting
# the page the same way.
# ---------------------------------------------------------------------------
$ws.PageSetup.Zoom = 88
